$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.164982199668884
$ws.Range("B1").Value = 2.420440673828125
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.37702488899231
$ws.Range("E1").Value = 1.235623598098755
